$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells keep their exact text formatting (avoid Excel
# auto-converting numeric-looking strings like "224.30" or "10.00" into numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "33.661.71"
$ws.Range("E2").Value = "  +6.85%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.772.47"
$ws.Range("E3").Value = "  +3.76%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.30"
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.555"
$ws.Range("E6").Value = "  +3.69%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.86"
$ws.Range("E8").Value = "  -0.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.64"
$ws.Range("E9").Value = "  +3.89%  "

$ws.Range("E10").Value = "  +2.53%  "

$ws.Range("E11").Value = "  +1.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0922"
$ws.Range("E12").Value = "  +1.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.024.15"
$ws.Range("E13").Value = "  +3.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.770.98"
$ws.Range("E14").Value = "  +3.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.622"
$ws.Range("E15").Value = "  +1.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.617.78"
$ws.Range("E16").Value = "  +6.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.00"
$ws.Range("E17").Value = "  -2.03%  "

$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.17"
$ws.Range("E19").Value = "  +1.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "249.18"
$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("E21").Value = "  +1.45%  "

$ws.Range("E22").Value = "  +0.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.21"
$ws.Range("E23").Value = "  +0.76%  "

$ws.Range("E24").Value = "  -2.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.13"
$ws.Range("E25").Value = "  -1.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.77"
$ws.Range("E26").Value = "  -0.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.33"
$ws.Range("E27").Value = "  +1.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("E28").Value = "  +0.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.91"
$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("E31").Value = "  -2.06%  "

$ws.Range("E32").Value = "  +1.98%  "

$ws.Range("E33").Value = "  +3.18%  "

$ws.Range("E34").Value = "  +4.45%  "

$ws.Range("E35").Value = "  +4.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.473.34"
$ws.Range("E36").Value = "  -3.23%  "

$ws.Range("E37").Value = "  +2.21%  "

$ws.Range("E38").Value = "  +2.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0184"
$ws.Range("E39").Value = "  +1.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.47"
$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.35"
$ws.Range("E41").Value = "  +1.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.68"
$ws.Range("E42").Value = "  -1.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.881"
$ws.Range("E43").Value = "  +3.28%  "

$ws.Range("E44").Value = "  +1.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.917.95"
$ws.Range("E47").Value = "  +3.67%  "

$ws.Range("E48").Value = "  +0.14%  "

$ws.Range("E49").Value = "  +1.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.54"
$ws.Range("E50").Value = "  +12.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.61"
$ws.Range("E51").Value = "  -2.90%  "

# Rows 45 and 46: Kaspa and WEMIXToken swapped ranking positions, with updated data
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0508"
$ws.Range("E45").Value = "  +0.75%  "

$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.07"
$ws.Range("E46").Value = "  +4.06%  "
